$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price figures stored as literal text (e.g. "67.278.44", "0.999", "2.66").
# Pre-format each touched D cell as Text so Excel keeps the literal digits/dots instead of
# re-parsing them as a number, matching the original inlineStr/text cell content & type.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

$ws.Range('D2').Value = '67.278.44'
$ws.Range('E2').Value = '  -4.82%  '
$ws.Range('D3').Value = '3.258.47'
$ws.Range('E3').Value = '  -7.43%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '599.44'
$ws.Range('E5').Value = '  -3.22%  '
$ws.Range('D6').Value = '151.21'
$ws.Range('E6').Value = '  -12.60%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.249.96'
$ws.Range('E8').Value = '  -7.56%  '
$ws.Range('E9').Value = '  -11.28%  '
$ws.Range('E10').Value = '  -13.25%  '
$ws.Range('D11').Value = '6.71'
$ws.Range('E11').Value = '  -4.92%  '
$ws.Range('E12').Value = '  -13.77%  '
$ws.Range('D13').Value = '38.13'
$ws.Range('E13').Value = '  -17.93%  '
$ws.Range('D14').Value = '0.0000242'
$ws.Range('E14').Value = '  -12.04%  '
$ws.Range('D15').Value = '3.779.87'
$ws.Range('E15').Value = '  -7.60%  '
$ws.Range('D16').Value = '67.287.50'
$ws.Range('E16').Value = '  -4.91%  '
$ws.Range('D17').Value = '3.258.54'
$ws.Range('E17').Value = '  -7.36%  '
$ws.Range('D18').Value = '542.52'
$ws.Range('E18').Value = '  -10.88%  '
$ws.Range('E19').Value = '  -6.06%  '
$ws.Range('D20').Value = '7.21'
$ws.Range('E20').Value = '  -13.91%  '
$ws.Range('E21').Value = '  -14.65%  '
$ws.Range('E22').Value = '  -13.73%  '
$ws.Range('D23').Value = '7.87'
$ws.Range('E23').Value = '  -14.29%  '
$ws.Range('E24').Value = '  -12.63%  '
$ws.Range('D25').Value = '13.50'
$ws.Range('E25').Value = '  -13.70%  '
$ws.Range('D27').Value = '3.27'
$ws.Range('E27').Value = '  -12.39%  '
$ws.Range('E28').Value = '  -11.02%  '
$ws.Range('D29').Value = '29.34'
$ws.Range('E29').Value = '  -12.92%  '
$ws.Range('E30').Value = '  -17.30%  '
$ws.Range('E31').Value = '  -10.85%  '
$ws.Range('E32').Value = '  -12.10%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '6.66'
$ws.Range('E33').Value = '  -17.72%  '
$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').Value = '543.43'
$ws.Range('E34').Value = '  -15.10%  '
$ws.Range('D35').Value = '5.72'
$ws.Range('E35').Value = '  -16.20%  '
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('D37').Value = '0.0445'
$ws.Range('E37').Value = '  -8.30%  '
$ws.Range('D38').Value = '53.20'
$ws.Range('E38').Value = '  -6.03%  '
$ws.Range('D39').Value = '0.0853'
$ws.Range('E39').Value = '  -14.41%  '
$ws.Range('E40').Value = '  -15.07%  '
$ws.Range('D41').Value = '0.129'
$ws.Range('E41').Value = '  -9.86%  '
$ws.Range('D42').Value = '2.925.96'
$ws.Range('E42').Value = '  -12.68%  '
$ws.Range('D43').Value = '2.66'
$ws.Range('E43').Value = '  -22.40%  '
$ws.Range('E44').Value = '  -16.02%  '
$ws.Range('D45').Value = '0.0₃0585'
$ws.Range('E45').Value = '  -18.44%  '
$ws.Range('D46').Value = '2.18'
$ws.Range('E46').Value = '  -14.01%  '
$ws.Range('D47').Value = '26.51'
$ws.Range('E47').Value = '  -16.62%  '
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('D49').Value = '127.95'
$ws.Range('E49').Value = '  -4.92%  '
$ws.Range('D50').Value = '2.35'
$ws.Range('E50').Value = '  -20.24%  '
$ws.Range('E51').Value = '  -12.68%  '
